$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.082.52'
$ws.Range('E2').Value = '  +1.55%  '
$ws.Range('D3').Value = '3.865.46'
$ws.Range('E3').Value = '  +1.13%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '474.23'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +10.97%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.29'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +10.45%  '
$ws.Range('E7').Value = '  +3.17%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.747'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.76%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.155'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000312'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -7.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.55'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.76%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.40'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.33%  '
$ws.Range('D14').Value = '4.480.98'
$ws.Range('E14').Value = '  +1.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.84'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.06%  '
$ws.Range('D16').Value = '3.888.17'
$ws.Range('E16').Value = '  +1.09%  '
$ws.Range('E17').Value = '  -0.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.13'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('E19').Value = '  +5.15%  '
$ws.Range('D20').Value = '67.278.84'
$ws.Range('E20').Value = '  +1.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '432.82'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.76%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.91'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.56%  '
$ws.Range('E23').Value = '  +6.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.66'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.87%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.61'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +9.22%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '38.07'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.96%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.40%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.96'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.56'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.45%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '732.39'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.91'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.135'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.79'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.59%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '43.55'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +12.06%  '
$ws.Range('E35').Value = '  +6.89%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '57.87'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.00%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.48'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0485'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.26%  '
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.349'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.91%  '
$ws.Range('B41').Value = 'ThetaToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.92'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.142'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.35%  '
$ws.Range('B43').Value = 'Fetch.AI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.57'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.97%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('B45').Value = 'PEPE'
$ws.Range('C45').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D45').Value = '0.0₃0672'
$ws.Range('E45').Value = '  -8.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.47'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.07%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.77'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.18%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.16'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.09%  '
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.20'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.79%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.89'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.78%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '144.03'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.01%  '
